# Add a new "chi-sq" worksheet with chi-squared statistical comparison
# tables (qPCR vs RDT, and ML model comparisons by "Types" grouping),
# positioned between "DNA - Test Matches" and "RDTqPCR Confustion Matrix".

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("DNA - Test Matches")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "chi-sq"

# ---- column widths (approximate character widths) ----
$ws.Columns.Item(1).ColumnWidth = 5.16
$ws.Columns.Item(2).ColumnWidth = 12.16
$ws.Columns.Item(3).ColumnWidth = 15.33
$ws.Columns.Item(4).ColumnWidth = 18.66
$ws.Columns.Item(5).ColumnWidth = 23.16
$ws.Columns.Item(6).ColumnWidth = 16.5
$ws.Columns.Item(7).ColumnWidth = 23.33
$ws.Columns.Item(8).ColumnWidth = 21.33
$ws.Columns.Item(9).ColumnWidth = 24.16

# ---- row 1: qPCR/RDT vs ML model headers ----
$ws.Range("B1").Value = "qPCR"
$ws.Range("C1").Value = "RDT"
$ws.Range("D1").Value = "Voting Ensemble`nTypes 2 and 7"
$ws.Range("E1").Value = "Extreme Random Trees`nTypes 2 and 7"
$ws.Range("F1").Value = "Voting Ensemble`nTypes 1 thru 24"
$ws.Range("G1").Value = "Extreme Random Trees`nTypes 1 thru 24"
$ws.Range("H1").Value = "Voting Ensemble`nTypes 3 and 5 and 10"
$ws.Range("I1").Value = "Extreme Random Trees`nTypes 3 and 5 and 10"

# ---- row 2: qPCR comparison row ----
$ws.Range("A2").Value = "qPCR"
$ws.Range("B2").Value = "67`n(2.715e-16)"
$ws.Range("C2").Value = "32.668`n(1.093e-08)"
$ws.Range("D2").Value = "7.1373`n(0.00755)"
$ws.Range("E2").Value = "7.1373`n(0.00755)"
$ws.Range("F2").Value = "8.963`n(0.002755)"
$ws.Range("G2").Value = "9.9844`n(0.001579)"
$ws.Range("H2").Value = "10.338`n(0.001303)"
$ws.Range("I2").Value = "10.866`n(0.0009797)"

# ---- row 3: RDT comparison row ----
$ws.Range("A3").Value = "RDT"
$ws.Range("B3").Value = "32.668`n(1.093e-08)"
$ws.Range("C3").Value = "67`n(2.715e-16)"
$ws.Range("D3").Value = "3.1485`n(0.076)"
$ws.Range("E3").Value = "3.1485`n(0.076)"
$ws.Range("F3").Value = "11.981`n(0.0005375)"
$ws.Range("G3").Value = "10.073`n(0.001505)"
$ws.Range("H3").Value = "12.55`n(0.0003962)"
$ws.Range("I3").Value = "15.5`n(8.249e-05)"

# ---- row 8-9: qPCR vs RDT mini header (merged) ----
$ws.Range("B8:C8").Merge()
$ws.Range("B8").Value = "qPCR"

$ws.Range("A9").Value = "RDT"
$ws.Range("B9:C9").Merge()
$ws.Range("B9").Value = "32.668`n(1.093e-08)"

# ---- rows 10-13: ML model (Voting Ensemble / Extreme Random Trees) by Types group ----
$ws.Range("A10").Value = "ML"
$ws.Range("B10").Value = "Voting Ensemble"
$ws.Range("C10").Value = "Extreme Random Trees"

$ws.Range("A11").Value = "Types 2 and 7"
$ws.Range("B11").Value = "7.1373`n(0.00755)"
$ws.Range("C11").Value = "7.1373`n(0.00755)"

$ws.Range("A12").Value = "Types 1 thru 24"
$ws.Range("B12").Value = "8.963`n(0.002755)"
$ws.Range("C12").Value = "9.9844`n(0.001579)"

$ws.Range("A13").Value = "Types 3 and 5 and 10"
$ws.Range("B13").Value = "10.338`n(0.001303)"
$ws.Range("C13").Value = "10.866`n(0.0009797)"

# row 14 trailing formatted (empty) cells
$ws.Range("A14").Value = ""
$ws.Range("C14").Value = ""

# ---- alignment: center everything used, wrap the larger text blocks ----
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4108

$ws.Range("D1:I3").HorizontalAlignment = -4108
$ws.Range("D1:I3").VerticalAlignment = -4108
$ws.Range("D1:I3").WrapText = $true

$ws.Range("A2:A3").HorizontalAlignment = -4108
$ws.Range("A2:A3").VerticalAlignment = -4108

$ws.Range("B2:I3").HorizontalAlignment = -4108
$ws.Range("B2:I3").VerticalAlignment = -4108
$ws.Range("B2:I3").WrapText = $true

$ws.Range("B8:C8").HorizontalAlignment = -4108
$ws.Range("B8:C8").VerticalAlignment = -4108

$ws.Range("A9").HorizontalAlignment = -4108
$ws.Range("A9").VerticalAlignment = -4108

$ws.Range("B9:C9").HorizontalAlignment = -4108
$ws.Range("B9:C9").VerticalAlignment = -4108
$ws.Range("B9:C9").WrapText = $true

$ws.Range("A10").HorizontalAlignment = -4108
$ws.Range("A10").VerticalAlignment = -4108

$ws.Range("B10:C10").HorizontalAlignment = -4108
$ws.Range("B10:C10").VerticalAlignment = -4108
$ws.Range("B10:C10").WrapText = $true

$ws.Range("A11:C13").HorizontalAlignment = -4108
$ws.Range("A11:C13").VerticalAlignment = -4108
$ws.Range("A11:C13").WrapText = $true

$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("A14").WrapText = $true

$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("C14").VerticalAlignment = -4108
$ws.Range("C14").WrapText = $true

# ---- row heights ----
$ws.Range("A1:I3").RowHeight = 32
$ws.Range("A9:I9").RowHeight = 32
$ws.Range("A10:I10").RowHeight = 32
$ws.Range("A11:I12").RowHeight = 64
$ws.Range("A13:I13").RowHeight = 96

# ---- view ----
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("A8:C13").Select()
